$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.698.32"
$ws.Range("E2").Value = "  +3.09%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.165.64"
$ws.Range("E3").Value = "  +3.54%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'538.96"
$ws.Range("E5").Value = "  +3.78%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'140.56"
$ws.Range("E6").Value = "  +4.43%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.13%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.512"
$ws.Range("E8").Value = "  +10.66%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +0.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.36%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.422"
$ws.Range("E11").Value = "  +5.34%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.73%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.706.67"
$ws.Range("E13").Value = "  +3.42%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'26.30"
$ws.Range("E14").Value = "  +5.34%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +7.23%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "58.747.45"
$ws.Range("E16").Value = "  +3.04%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "'6.26"
$ws.Range("E17").Value = "  +7.39%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.152.05"
$ws.Range("E18").Value = "  +3.24%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'13.09"
$ws.Range("E19").Value = "  +6.12%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'8.22"
$ws.Range("E20").Value = "  +6.41%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'378.42"
$ws.Range("E21").Value = "  +9.11%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.16%  "

# Row 23 - LEO
$ws.Range("D23").Value = "'5.76"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'70.56"
$ws.Range("E24").Value = "  +2.48%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +4.81%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +3.23%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.32%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'8.10"
$ws.Range("E28").Value = "  +14.18%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +5.91%  "

# Row 30 - RenderToken (was PancakeSwap)
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'6.22"
$ws.Range("E30").Value = "  +8.83%  "

# Row 31 - PancakeSwap (was RenderToken)
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.90"
$ws.Range("E31").Value = "  +3.31%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'22.00"
$ws.Range("E32").Value = "  +5.82%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "'5.18"
$ws.Range("E33").Value = "  +8.75%  "

# Row 34 - Fetch.AI
$ws.Range("D34").Value = "'1.19"
$ws.Range("E34").Value = "  +7.22%  "

# Row 35 - Monero
$ws.Range("D35").Value = "'161.68"
$ws.Range("E35").Value = "  +2.37%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +5.67%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'1.37"
$ws.Range("E37").Value = "  +12.63%  "

# Row 38 - EnergySwap
$ws.Range("D38").Value = "'25.60"
$ws.Range("E38").Value = "  +2.14%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +8.89%  "

# Row 40 - Maker
$ws.Range("D40").Value = "2.658.43"
$ws.Range("E40").Value = "  +11.16%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "'0.0684"
$ws.Range("E41").Value = "  +5.20%  "

# Row 42 - Filecoin
$ws.Range("D42").Value = "'4.23"
$ws.Range("E42").Value = "  +6.14%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'38.70"
$ws.Range("E43").Value = "  +6.41%  "

# Row 44 - Mantle
$ws.Range("D44").Value = "'0.709"
$ws.Range("E44").Value = "  +3.46%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +7.15%  "

# Row 46 - FirstDigitalUSD
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.10%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "'0.104"
$ws.Range("E47").Value = "  +13.97%  "

# Row 48 - Cosmos
$ws.Range("E48").Value = "  +5.02%  "

# Row 49 - ONDO
$ws.Range("D49").Value = "'0.986"
$ws.Range("E49").Value = "  +6.70%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'20.31"
$ws.Range("E50").Value = "  +5.67%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "'0.757"
$ws.Range("E51").Value = "  +2.84%  "
